# Auto-generated edit script: updates cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.127.29'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.87'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.81'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  +7.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3734'
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07327'
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8612'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.00'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.827.79'
$ws.Range("E12").Value = '  -2.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.699'
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.96'
$ws.Range("E14").Value = '  +5.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.345'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07081'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008842'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.128.59'
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.194'
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.03'
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.999'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.68'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.224'
$ws.Range("E26").Value = '  +5.65%  '
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.270'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.36'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08873'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7719'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.973'
$ws.Range("E33").Value = '  +6.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.467'
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01966'
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05291'
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5366'
$ws.Range("E39").Value = '  +7.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.195'
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1718'
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5234'
$ws.Range("E43").Value = '  +11.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.630'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.70'
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.001'
$ws.Range("E46").Value = '  +11.52%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06509'
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.679'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9999'
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9229'
$ws.Range("E51").Value = '  +1.36%  '
